$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$n = 226
$vals = New-Object 'object[,]' $n,1
for ($i = 0; $i -lt $n; $i++) {
    $vals[$i,0] = $i
}
$ws.Range("U2:U227").Value = $vals

$ws.Range("U2").Select()
